$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:E and G (rows 2-8). F (Win) column is unchanged.
$data = @{
    2 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    3 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    4 = @(1.459612070389937, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 2.42670696938877)
    5 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 13.71653804550039)
    6 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 14.36450238910742)
    7 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    8 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 6.740334628841572)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}
